$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Shift-Cell($row, $fromCol, $toCol) {
    # Capture the source cell's current value, grab its format onto the
    # clipboard, write the value into the target, then stamp the target
    # with the copied format as the very last step (setting .Value2 can
    # otherwise cause the target's own style to be re-normalised).
    $val = $ws.Cells.Item($row, $fromCol).Value2
    $ws.Cells.Item($row, $fromCol).Copy() | Out-Null
    $ws.Cells.Item($row, $toCol).Value2 = $val
    $ws.Cells.Item($row, $toCol).PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false
}

function Set-Cell($row, $col, $value, $formatFromCol) {
    $ws.Cells.Item($row, $formatFromCol).Copy() | Out-Null
    $ws.Cells.Item($row, $col).Value2 = $value
    $ws.Cells.Item($row, $col).PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false
}

for ($row = 1; $row -le 5; $row++) {
    # Shift existing group columns one place to the right: D->E, C->D, B->C
    Shift-Cell $row 4 5
    Shift-Cell $row 3 4
    Shift-Cell $row 2 3
}

# New column B: header "OTROS" (same style as the other headers), and
# "---" (no access) for every data row, matching the other data cells.
# (Format source columns are picked from cells whose style is already in
# its final resting state at this point in the script.)
Set-Cell 1 2 "OTROS" 3
Set-Cell 2 2 "---" 4
Set-Cell 3 2 "---" 3
Set-Cell 4 2 "---" 3
Set-Cell 5 2 "---" 3

# Match the saved selection state.
$ws.Range("D1").Select() | Out-Null
